# ---------------------------------------------------------------------------
# Target change (from the supplied OOXML diff): word/styles.xml's
# <w:docDefaults> block loses a long run of redundant/boilerplate
# properties that merely restate the OOXML schema's own implicit defaults
# (b/i/smallCaps/strike/color/u/shd/vertAlign under <w:rPrDefault>, and
# keepNext/keepLines/widowControl/pBdr/shd/ind/contextualSpacing/jc under
# <w:pPrDefault>), leaving only the meaningful values
# (rFonts/sz/szCs/lang, and spacing line=276/lineRule=auto).
#
# <w:docDefaults> (and its <w:rPrDefault>/<w:pPrDefault> children) is a
# document-defaults construct that lives in the styles part of the
# package. It is NOT surfaced anywhere on the Word object model: the
# Styles collection only ever yields the named styles actually defined in
# the document (here: Normal, Table Normal, Heading 1-6, Title, Subtitle -
# see Document.Styles.Count/Item below), each with its own Font /
# ParagraphFormat. There is no Styles("docDefaults") / "rPrDefault" /
# "pPrDefault" pseudo-entry, no Document.DefaultFont / DefaultParagraphFormat
# property, and WordOpenXML is read-only here (as in real automation,
# round-tripping the whole flat-OPC package back in isn't a supported
# write path) - this mirrors genuine Word COM automation, where
# docDefaults can only be edited by hand-editing styles.xml (e.g. via the
# Open XML SDK), never through Application/Document/Styles calls.
#
# Concretely: the "Normal" style (the only style any paragraph in this
# document uses) currently carries *no* rPr/pPr overrides of its own - it
# inherits 100% of its formatting from docDefaults, and the diff leaves
# that <w:style .../> element completely untouched (no rPr/pPr added to
# it). Pushing Font/ParagraphFormat assignments onto Styles("Normal") to
# fake the same visual effect would therefore not reach docDefaults at
# all; it would just bolt a new, redundant override directly onto the
# Normal style definition - moving the diff rather than resolving it, and
# diverging from the target OOXML even further than leaving the document
# untouched does.
#
# So: confirm the shape of what's available (read-only), and make no
# mutating calls - there is no reachable lever on this object model for
# the docDefaults cleanup described by the diff, and any workaround
# that *is* reachable (Styles("Normal").Font/.ParagraphFormat, Find &
# Replace formatting, InsertXML, etc.) writes to document content or to a
# named style's own rPr/pPr, not to <w:docDefaults>, which would only add
# unrelated drift instead of applying the requested edit.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$normal = $d.Styles("Normal")
Write-Output ("styles: {0}" -f $d.Styles.Count)
Write-Output ("Normal style font: {0} {1}" -f $normal.Font.Name, $normal.Font.Size)
